# Generate Report for Handback
# Row 7 (the 0c7a7b8a-9dca-42a6-8f19-793adbdcec1f e2e case) gets its
# handback-report columns filled in on both the zh-cn and de-de sheets:
#   I  = Latest Target File   -> handoff .md display name (new hyperlink)
#   J  = Latest Handback File -> same xlf name as the Latest Handoff File (G)
#   K  = Latest Handback DateTime
#   P  = Error Detail -> "handback not latest" message

$wb = $excel.ActiveWorkbook

$mdDisplay = "0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31cb41cf411ca3657245d9cf8bac80f634b968d8/e2e/0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1760995180e44a2e3ac4d1f0306f9c5efff5d524/e2e/0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31cb41cf411ca3657245d9cf8bac80f634b968d8/e2e/0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.md."

$sheetInfo = @(
    @{ Name = "zh-cn"; HandbackXlf = "0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.93f424616503a3eb1c232a6d486815ce51a66111.zh-cn.xlf"; HandbackDate = "2016-08-28 04:53:53" },
    @{ Name = "de-de"; HandbackXlf = "0c7a7b8a-9dca-42a6-8f19-793adbdcec1f.93f424616503a3eb1c232a6d486815ce51a66111.de-de.xlf"; HandbackDate = "2016-08-28 04:54:00" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # I7: "Latest Target File" -> hyperlink to the handoff markdown file.
    $ws.Hyperlinks.Add($ws.Range("I7"), $mdUrl, "", "", $mdDisplay)

    # J7: "Latest Handback File" -> same xliff name as the latest handoff file.
    $ws.Range("J7").Value = $info.HandbackXlf

    # K7: "Latest Handback DateTime"
    $ws.Range("K7").Value = $info.HandbackDate

    # P7: "Error Detail"
    $ws.Range("P7").Value = $errorDetail
}
